# TimeSheet.xlsx - add a new (blank) activity row to the last month block
# (Bahman 98, rows 83-88) just above the "Total Hours" row, update the
# Python Tracker Connection hours, and refresh the Paid/Not-Paid totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "* Python Tracker Connection" hours: 2 -> 4
$ws.Range("C84").Value = 4

# Insert a new blank activity row above the old "Total Hours" row (row 86).
# This pushes the Total Hours / Paid / Not Paid rows down by one (86->87,
# 87->88, 88->89), matching how Excel shifts rows on insert.
$ws.Rows(86).Insert()

# The newly inserted row is blank by default but needs the same borders /
# alignment as the other activity rows (B: left align, C: center align) -
# copy that formatting down from the row above (row 85).
$ws.Range("B85:C85").Copy()
$ws.Range("B86:C86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Total Hours" row is now row 87 - widen its SUM to cover the new row.
$ws.Range("C87").Formula = "=SUM(C84:C86)"

# "Not Paid" (@Home) row is now row 89 - hours: 2 -> 4
$ws.Range("D89").Value = 4

# Keep the on-screen view in sync with the new layout (scrolled to the
# last block, cursor resting one row higher than before since the block
# grew by a row).
$ws.Range("E87").Select()
